# Início da implementação da funcionalidade de leitura da revista de patente
#
# Adds 5 new "Data / Quantidade de horas" rows (19-23) to Planilha1, right
# after the existing data (which ended at row 18), and updates the sheet
# selection to match where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: Excel date serial numbers in column A, hour fractions (of a day)
# in column B - same shape as all the existing rows above them.
$dates = @(41575, 41577, 41578, 41580, 41581)
$hours = @(0.05555555555555555, 0.027777777777777776, 0.020833333333333332, 0.06944444444444443, 0.09027777777777778)

$startRow = 19
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $hours[$i]

    # Match the date / time number formats used by the rows above (A: date,
    # B: elapsed hours:minutes).
    $ws.Cells.Item($row, 1).NumberFormat = "m/d/yy"
    $ws.Cells.Item($row, 2).NumberFormat = "h:mm"
}

# The selection now spans the whole "Quantidade de horas" column down to the
# newly-added last row.
$ws.Range("B2:B23").Select()
